# Update "想去人数" (want-to-go count) figures across sheets to reflect
# the latest generated output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 1400
$wsExpo.Range("F9").Value  = 1844
$wsExpo.Range("F10").Value = 490
$wsExpo.Range("F12").Value = 22
$wsExpo.Range("F16").Value = 7028
$wsExpo.Range("F19").Value = 174
$wsExpo.Range("F21").Value = 1718
$wsExpo.Range("F28").Value = 329

# Sheet "演出" (Performance)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 354

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value  = 354
$wsAll.Range("F10").Value = 1400
$wsAll.Range("F14").Value = 1844
$wsAll.Range("F17").Value = 490
$wsAll.Range("F19").Value = 22
$wsAll.Range("F24").Value = 7029
$wsAll.Range("F29").Value = 1718
$wsAll.Range("F36").Value = 329
